# Project2/ProductOrders.xlsx edit
#
# Fixes a mislabeled product name on the "Orders" sheet (row for the 4th
# order had been using the typo'd "Ipoh Coff" string instead of the
# correct "Ipoh Coffee" product name already used elsewhere), and starts
# populating the Status column with the add-to-cart outcome ("Succeeded"
# / "Failed") now that the site-iteration script can detect 404s and
# out-of-stock items.

$wb = $excel.ActiveWorkbook

$orders = $wb.Worksheets.Item("Orders")

# Row 9 ("Ipoh Coff") should reference the same product as row 3 ("Ipoh Coffee").
$orders.Range("A9").Value = "Ipoh Coffee"

# Record add-to-cart results for the first couple of rows.
$orders.Range("C2").Value = "Succeeded"
$orders.Range("C3").Value = "Failed"

# Bring the Orders tab back into focus (it had been left on Address).
$orders.Activate()
